$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- E8: set Acceptance to "Accepted" ---
$ws.Range("E8").Value = "Accepted"

# --- G8: append new comment text with mixed bold requirement IDs ---
$fullText = "Added 4 new requirements Req_PO1_DGC_SRS_014_V01, Req_PO1_DGC_SRS_015_V01, Req_PO1_DGC_SRS_016_V01 and Req_PO1_DGC_SRS_017_V01 to indicate the ranges of the horizontal and vertical cursor positions and what should happen if they were exceeded`n`nMali 19/2/2020: I didn't mean that, I meant for example Req_PO1_DGC_SRS_009_V01 the requirement shall mention the value of x and y`nMina 20/02/2020: The values cannot be mentioned exactly/explicitly since they're a range. But the requirements Req_PO1_DGC_SRS_014_V01 and Req_PO1_DGC_SRS_016_V01 were updated to be more descriptive."

$cell = $ws.Range("G8")
$cell.Value = $fullText

$boldRanges = @(
    @(26, 23),
    @(51, 23),
    @(76, 23),
    @(104, 23),
    @(487, 23),
    @(515, 23)
)

foreach ($r in $boldRanges) {
    $start = $r[0]
    $len = $r[1]
    $chars = $cell.Characters($start, $len)
    $chars.Font.Bold = $true
}

# --- sheet view: scroll/selection update ---
$ws.Range("G8").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 6

# --- row heights (re-flowed by Excel's autofit after the edits) ---
$ws.Rows.Item(2).RowHeight = 28.8
$ws.Rows.Item(3).RowHeight = 57.6
$ws.Rows.Item(4).RowHeight = 43.2
$ws.Rows.Item(5).RowHeight = 43.2
$ws.Rows.Item(6).RowHeight = 57.6
$ws.Rows.Item(7).RowHeight = 57.6
$ws.Rows.Item(8).RowHeight = 244.8

Write-Output "done"
